$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 4702.5
$ws.Range("I7").Value = 4405
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 4405
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -4293
$ws.Range("N7").Value = -5224
$ws.Range("H14").Value = 4702.5
$ws.Range("I14").Value = 4405
$ws.Range("J14").Value = 5000
$ws.Range("K14").Value = 4405
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = -4214
$ws.Range("N14").Value = -5382
$ws.Range("H40").Value = 1578.2354
$ws.Range("I40").Value = 1158.7059
$ws.Range("J40").Value = 1997.7646
$ws.Range("K40").Value = 1158.7059
$ws.Range("L40").Value = 1997.7646
$ws.Range("M40").Value = -983.7058999999999
$ws.Range("N40").Value = -2347.7646
$ws.Range("H46").Value = 1293.25
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1293.25
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3879.75
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -4117.75
$ws.Range("H60").Value = 1293.25
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 1293.25
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 3879.75
$ws.Range("M60").Value = $null
$ws.Range("N60").Value = -4847.75
$ws.Range("H62").Value = 13891629
$ws.Range("I62").Value = 16669515
$ws.Range("J62").Value = 2200
$ws.Range("K62").Value = 16669515
$ws.Range("L62").Value = 2200
$ws.Range("M62").Value = -16668891
$ws.Range("N62").Value = -3448
$ws.Range("H64").Value = 3274.1614
$ws.Range("I64").Value = 3131.5264
$ws.Range("K64").Value = 3131.5264
$ws.Range("M64").Value = -2883.5264
$ws.Range("H65").Value = 13891629
$ws.Range("I65").Value = 16669515
$ws.Range("J65").Value = 2200
$ws.Range("K65").Value = 83347575
$ws.Range("L65").Value = 11000
$ws.Range("M65").Value = -83344455
$ws.Range("N65").Value = -17240
$ws.Range("H67").Value = 3274.1614
$ws.Range("I67").Value = 3131.5264
$ws.Range("K67").Value = 3131.5264
$ws.Range("M67").Value = -2273.5264
$ws.Range("H74").Value = 7387.5
$ws.Range("I74").Value = 3050
$ws.Range("J74").Value = 8833.333000000001
$ws.Range("K74").Value = 3050
$ws.Range("L74").Value = 8833.333000000001
$ws.Range("M74").Value = -2114
$ws.Range("N74").Value = -10705.333
$ws.Range("H77").Value = 7387.5
$ws.Range("I77").Value = 3050
$ws.Range("J77").Value = 8833.333000000001
$ws.Range("K77").Value = 15250
$ws.Range("L77").Value = 44166.665
$ws.Range("M77").Value = -10570
$ws.Range("N77").Value = -53526.665
$ws.Range("H100").Value = 8243
$ws.Range("I100").Value = 13443.125
$ws.Range("J100").Value = 2300
$ws.Range("K100").Value = 13443.125
$ws.Range("L100").Value = 2300
$ws.Range("M100").Value = -12902.125
$ws.Range("N100").Value = -3382

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2667.1333
$ws.Range("I26").Value = 1503.5
$ws.Range("J26").Value = 2846.1538
$ws.Range("K26").Value = 1503.5
$ws.Range("L26").Value = 2846.1538
$ws.Range("M26").Value = -1173.5
$ws.Range("N26").Value = -3506.1538
$ws.Range("H61").Value = 2138.32
$ws.Range("I61").Value = 2284.923
$ws.Range("J61").Value = 1979.5
$ws.Range("K61").Value = 2284.923
$ws.Range("L61").Value = 1979.5
$ws.Range("M61").Value = -2072.923
$ws.Range("N61").Value = -2403.5
$ws.Range("H110").Value = 862.875
$ws.Range("I110").Value = 509.9565
$ws.Range("J110").Value = 8980
$ws.Range("K110").Value = 509.9565
$ws.Range("L110").Value = 8980
$ws.Range("M110").Value = 1535.0435
$ws.Range("N110").Value = -13070
$ws.Range("H122").Value = 1431.7142
$ws.Range("I122").Value = 1012
$ws.Range("J122").Value = 1599.6
$ws.Range("K122").Value = 3036
$ws.Range("L122").Value = 4798.799999999999
$ws.Range("M122").Value = -586
$ws.Range("N122").Value = -9698.799999999999
$ws.Range("H132").Value = 6557.0435
$ws.Range("I132").Value = 1801.5
$ws.Range("J132").Value = 9093.333000000001
$ws.Range("K132").Value = 5404.5
$ws.Range("L132").Value = 27279.999
$ws.Range("M132").Value = -2874.5
$ws.Range("N132").Value = -32339.999
$ws.Range("H136").Value = 2138.32
$ws.Range("I136").Value = 2284.923
$ws.Range("J136").Value = 1979.5
$ws.Range("K136").Value = 6854.768999999999
$ws.Range("L136").Value = 5938.5
$ws.Range("M136").Value = -4304.768999999999
$ws.Range("N136").Value = -11038.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 260.1875
$ws.Range("I64").Value = 309.125
$ws.Range("J64").Value = 211.25
$ws.Range("K64").Value = 309.125
$ws.Range("L64").Value = 211.25
$ws.Range("M64").Value = -84.125
$ws.Range("N64").Value = -661.25
$ws.Range("H67").Value = 260.1875
$ws.Range("I67").Value = 309.125
$ws.Range("J67").Value = 211.25
$ws.Range("K67").Value = 309.125
$ws.Range("L67").Value = 211.25
$ws.Range("M67").Value = 470.875
$ws.Range("N67").Value = -1771.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2750
$ws.Range("I10").Value = 2750
$ws.Range("K10").Value = 2750
$ws.Range("M10").Value = -2611
$ws.Range("H14").Value = 8888.429
$ws.Range("J14").Value = 8888.429
$ws.Range("L14").Value = 8888.429
$ws.Range("N14").Value = -9228.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1000
$ws.Range("J49").Value = 1000
$ws.Range("L49").Value = 3000
$ws.Range("N49").Value = -3312
$ws.Range("H58").Value = 26896.154
$ws.Range("J58").Value = 29108.334
$ws.Range("L58").Value = 87325.00199999999
$ws.Range("N58").Value = -87581.00199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17247454
$ws.Range("I122").Value = 22734046
$ws.Range("J122").Value = 3884.4285
$ws.Range("K122").Value = 68202138
$ws.Range("L122").Value = 11653.2855
$ws.Range("M122").Value = -68199688
$ws.Range("N122").Value = -16553.2855
$ws.Range("H132").Value = 5433.607
$ws.Range("I132").Value = 1495
$ws.Range("J132").Value = 11520.546
$ws.Range("K132").Value = 4485
$ws.Range("L132").Value = 34561.638
$ws.Range("M132").Value = -1955
$ws.Range("N132").Value = -39621.638

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1487.4706
$ws.Range("I16").Value = 913.3333
$ws.Range("K16").Value = 913.3333
$ws.Range("M16").Value = -743.3333
$ws.Range("H136").Value = 3449.82
$ws.Range("I136").Value = 3079.2432
$ws.Range("K136").Value = 9237.729599999999
$ws.Range("M136").Value = -6687.729599999999
